$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: ticket/schedule id changes from 1000001 to 1001001
$ws.Range("A4").Value = 1001001

# New row 5 - duplicate the formatting of row 4 (copy+insert keeps the
# existing date style) then overwrite the cells that differ.
$ws.Rows("4:4").Copy() | Out-Null
$ws.Rows("5:5").Insert(-4121) | Out-Null
$ws.Range("A5").Value = 1002001
$ws.Range("B5").Value = "GACHA"
$ws.Range("C5").Value = 45658
$ws.Range("D5").Value = 45658
$ws.Range("E5").Value = 45736.999988425923
$ws.Range("F5").Value = 45736.999988425923

# New row 6 - same pattern as row 5
$ws.Rows("4:4").Copy() | Out-Null
$ws.Rows("6:6").Insert(-4121) | Out-Null
$ws.Range("A6").Value = 1002002
$ws.Range("B6").Value = "GACHA"
$ws.Range("C6").Value = 45658
$ws.Range("D6").Value = 45658
$ws.Range("E6").Value = 45736.999988425923
$ws.Range("F6").Value = 45736.999988425923

# Move the active selection like the saved workbook shows
$ws.Range("I7").Select() | Out-Null
